$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 984.1642000000001
$ws.Range("J17").Value = 831.4545000000001
$ws.Range("L17").Value = 2494.3635
$ws.Range("N17").Value = -2830.3635
$ws.Range("H62").Value = 2425
$ws.Range("I62").Value = 1994.75
$ws.Range("J62").Value = 2998.6667
$ws.Range("K62").Value = 1994.75
$ws.Range("L62").Value = 2998.6667
$ws.Range("M62").Value = -1370.75
$ws.Range("N62").Value = -4246.6667
$ws.Range("H64").Value = 3050
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3100
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3100
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3596
$ws.Range("H65").Value = 2425
$ws.Range("I65").Value = 1994.75
$ws.Range("J65").Value = 2998.6667
$ws.Range("K65").Value = 9973.75
$ws.Range("L65").Value = 14993.3335
$ws.Range("M65").Value = -6853.75
$ws.Range("N65").Value = -21233.3335
$ws.Range("H67").Value = 3050
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3100
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3100
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -4816
$ws.Range("H70").Value = 12028.111
$ws.Range("I70").Value = 699.6667
$ws.Range("J70").Value = 17692.334
$ws.Range("K70").Value = 2099.0001
$ws.Range("L70").Value = 53077.00199999999
$ws.Range("M70").Value = -1829.0001
$ws.Range("N70").Value = -53617.00199999999
$ws.Range("H73").Value = 12028.111
$ws.Range("I73").Value = 699.6667
$ws.Range("J73").Value = 17692.334
$ws.Range("K73").Value = 2099.0001
$ws.Range("L73").Value = 53077.00199999999
$ws.Range("M73").Value = -1163.0001
$ws.Range("N73").Value = -54949.00199999999
$ws.Range("H100").Value = 1761.9
$ws.Range("I100").Value = 1624.3334
$ws.Range("K100").Value = 1624.3334
$ws.Range("M100").Value = -1083.3334
$ws.Range("H113").Value = 36125.5
$ws.Range("I113").Value = 36125.5
$ws.Range("K113").Value = 36125.5
$ws.Range("M113").Value = -32871.5
$ws.Range("H123").Value = 43695
$ws.Range("J123").Value = 43695
$ws.Range("L123").Value = 43695
$ws.Range("N123").Value = -53495
$ws.Range("H132").Value = 1202.6774
$ws.Range("I132").Value = 1088.6786
$ws.Range("K132").Value = 3266.0358
$ws.Range("M132").Value = -736.0357999999997
$ws.Range("H138").Value = 2957.8108
$ws.Range("J138").Value = 3867.8
$ws.Range("L138").Value = 11603.4
$ws.Range("N138").Value = -21883.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1661813.8
$ws.Range("I2").Value = 2326159.8
$ws.Range("K2").Value = 2326159.8
$ws.Range("M2").Value = -2326046.8
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 37
$ws.Range("K4").Value = 37
$ws.Range("M4").Value = 79
$ws.Range("H8").Value = 20000000
$ws.Range("I8").Value = 20000000
$ws.Range("K8").Value = 20000000
$ws.Range("M8").Value = -19999856
$ws.Range("H11").Value = 10000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H13").Value = 23333334
$ws.Range("I13").Value = 23333334
$ws.Range("K13").Value = 23333334
$ws.Range("M13").Value = -23333190
$ws.Range("H61").Value = 3442.2
$ws.Range("J61").Value = 6000
$ws.Range("L61").Value = 6000
$ws.Range("N61").Value = -6424
$ws.Range("H74").Value = 1758.1282
$ws.Range("I74").Value = 1688
$ws.Range("J74").Value = 2599.6667
$ws.Range("K74").Value = 1688
$ws.Range("L74").Value = 2599.6667
$ws.Range("M74").Value = -814
$ws.Range("N74").Value = -4347.6667
$ws.Range("H77").Value = 1758.1282
$ws.Range("I77").Value = 1688
$ws.Range("J77").Value = 2599.6667
$ws.Range("K77").Value = 8440
$ws.Range("L77").Value = 12998.3335
$ws.Range("M77").Value = -4072
$ws.Range("N77").Value = -21734.3335
$ws.Range("H116").Value = 1661813.8
$ws.Range("I116").Value = 2326159.8
$ws.Range("K116").Value = 2326159.8
$ws.Range("M116").Value = -2323865.8
$ws.Range("H132").Value = 1864.579
$ws.Range("J132").Value = 2914
$ws.Range("L132").Value = 8742
$ws.Range("N132").Value = -13802
$ws.Range("H136").Value = 3442.2
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1661813.8
$ws.Range("I3").Value = 2326159.8
$ws.Range("K3").Value = 2326159.8
$ws.Range("M3").Value = -2326045.8
$ws.Range("H107").Value = 2519.4
$ws.Range("I107").Value = 2738.4285
$ws.Range("J107").Value = 2008.3334
$ws.Range("K107").Value = 2738.4285
$ws.Range("L107").Value = 2008.3334
$ws.Range("M107").Value = -818.4285
$ws.Range("N107").Value = -5848.3334
$ws.Range("H134").Value = 16312.667
$ws.Range("I134").Value = 16476.75
$ws.Range("K134").Value = 49430.25
$ws.Range("M134").Value = -46895.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 350
$ws.Range("I7").Value = 175
$ws.Range("J7").Value = 525
$ws.Range("K7").Value = 175
$ws.Range("L7").Value = 525
$ws.Range("M7").Value = -62
$ws.Range("N7").Value = -751
$ws.Range("H16").Value = 718.0769
$ws.Range("I16").Value = 743.625
$ws.Range("K16").Value = 743.625
$ws.Range("M16").Value = -456.625
$ws.Range("H58").Value = 3624767.2
$ws.Range("J58").Value = 2557
$ws.Range("L58").Value = 2557
$ws.Range("N58").Value = -2963
$ws.Range("H105").Value = 1655.5
$ws.Range("I105").Value = 1300
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 1300
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = 447
$ws.Range("N105").Value = -5505
$ws.Range("H107").Value = 324.33334
$ws.Range("I107").Value = 333.78946
$ws.Range("K107").Value = 333.78946
$ws.Range("M107").Value = 1586.21054
$ws.Range("H113").Value = 718.0769
$ws.Range("I113").Value = 743.625
$ws.Range("K113").Value = 743.625
$ws.Range("M113").Value = 1426.375
$ws.Range("H132").Value = 2145.9644
$ws.Range("I132").Value = 1167.1875
$ws.Range("J132").Value = 3451
$ws.Range("K132").Value = 3501.5625
$ws.Range("L132").Value = 10353
$ws.Range("M132").Value = -971.5625
$ws.Range("N132").Value = -15413
$ws.Range("H134").Value = 998.3333
$ws.Range("I134").Value = 998.3333
$ws.Range("K134").Value = 2994.9999
$ws.Range("M134").Value = -459.9998999999998
$ws.Range("H136").Value = 3624767.2
$ws.Range("J136").Value = 2557
$ws.Range("L136").Value = 7671
$ws.Range("N136").Value = -12771

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 753.11
$ws.Range("J131").Value = 779.3077
$ws.Range("L131").Value = 2337.9231
$ws.Range("N131").Value = -12417.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3849701
$ws.Range("I132").Value = 3849701
$ws.Range("K132").Value = 11549103
$ws.Range("M132").Value = -11546573

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 3380002.8
$ws.Range("I10").Value = 10000000
$ws.Range("J10").Value = 70004
$ws.Range("K10").Value = 10000000
$ws.Range("L10").Value = 70004
$ws.Range("M10").Value = -9999860
$ws.Range("N10").Value = -70284
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H132").Value = 2224
$ws.Range("I132").Value = 2403.7778
$ws.Range("J132").Value = 2116.1333
$ws.Range("K132").Value = 7211.3334
$ws.Range("L132").Value = 6348.3999
$ws.Range("M132").Value = -4681.3334
$ws.Range("N132").Value = -11408.3999
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("H136").Value = 3631.9
$ws.Range("I136").Value = 2136.6667
$ws.Range("K136").Value = 6410.000100000001
$ws.Range("M136").Value = -3860.000100000001
$ws.Range("H139").Value = 42500
$ws.Range("H141").Value = 30650
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws.Range("N133").Value = -94060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2050002.4
$ws.Range("J3").Value = 62503
$ws.Range("L3").Value = 62503
$ws.Range("N3").Value = -62731
$ws.Range("H13").Value = 3398
$ws.Range("I13").Value = 199
$ws.Range("J13").Value = 4997.5
$ws.Range("K13").Value = 199
$ws.Range("L13").Value = 4997.5
$ws.Range("M13").Value = -59
$ws.Range("H122").Value = 98433.875
$ws.Range("I122").Value = 130494.336
$ws.Range("J122").Value = 2252.5
$ws.Range("K122").Value = 391483.008
$ws.Range("L122").Value = 6757.5
$ws.Range("M122").Value = -389033.008
$ws.Range("N122").Value = -11657.5
$ws.Range("H132").Value = 1418.8695
$ws.Range("I132").Value = 924.2778
$ws.Range("J132").Value = 3199.4
$ws.Range("K132").Value = 2772.8334
$ws.Range("L132").Value = 9598.200000000001
$ws.Range("M132").Value = -242.8334
$ws.Range("N132").Value = -14658.2
$ws.Range("H136").Value = 3433.5881
$ws.Range("I136").Value = 4310.375
$ws.Range("J136").Value = 2654.2222
$ws.Range("K136").Value = 12931.125
$ws.Range("L136").Value = 7962.6666
$ws.Range("M136").Value = -10381.125
$ws.Range("N136").Value = -13062.6666
$ws.Range("N13").Value = -5277.5
